# Inserts a new weekly price record for "Haba" (Vega Modelo de Temuco) as
# row 40, pushing the existing rows 40-52 down to 41-53 (the sheet keeps a
# reverse-chronological / most-recent-first ordering, so the newest record
# is inserted near the top of the data block rather than appended at the
# bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40..52 down to 41..53, leaving a blank row 40 to fill in.
$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44518
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 125
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = 8000
$ws.Range("N40").Value = "$/saco 25 kilos"
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 320
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
